$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "300.60"
Set-TextValue $ws.Range("E2") "-4.56%"
Set-TextValue $ws.Range("G2") "10"

# Row 3
Set-TextValue $ws.Range("D3") "35.06"
Set-TextValue $ws.Range("E3") "-1.61%"
Set-TextValue $ws.Range("G3") "10"

# Row 4
Set-TextValue $ws.Range("D4") "5.046"
Set-TextValue $ws.Range("E4") "-1.88%"
Set-TextValue $ws.Range("G4") "10"

# Row 5
Set-TextValue $ws.Range("D5") "0.07940"
Set-TextValue $ws.Range("E5") "-2.09%"
Set-TextValue $ws.Range("G5") "10"

# Row 6
Set-TextValue $ws.Range("E6") "-10.58%"
Set-TextValue $ws.Range("G6") "10"

# Row 7
Set-TextValue $ws.Range("B7") "KuCoinToken"
Set-TextValue $ws.Range("C7") "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue $ws.Range("D7") "7.726"
Set-TextValue $ws.Range("E7") "-3.61%"
Set-TextValue $ws.Range("G7") "10"

# Row 8
Set-TextValue $ws.Range("B8") "BTSEToken"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D8") "2.950"
Set-TextValue $ws.Range("E8") "6.05%"
Set-TextValue $ws.Range("G8") "10"

# Row 9
Set-TextValue $ws.Range("B9") "MXToken"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D9") "0.9239"
Set-TextValue $ws.Range("E9") "-0.32%"
Set-TextValue $ws.Range("G9") "10"

# Row 10
Set-TextValue $ws.Range("B10") "LiechtensteinCryptoassetsExchange"
Set-TextValue $ws.Range("C10") "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1316"
Set-TextValue $ws.Range("E10") "28.26%"
Set-TextValue $ws.Range("G10") "10"

# Row 11
Set-TextValue $ws.Range("B11") "WazirX"
Set-TextValue $ws.Range("C11") "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1840"
Set-TextValue $ws.Range("E11") "-2.34%"
Set-TextValue $ws.Range("G11") "10"

# Row 12
Set-TextValue $ws.Range("B12") "MandalaExchangeToken"
Set-TextValue $ws.Range("C12") "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.09618"
Set-TextValue $ws.Range("E12") "3.69%"
Set-TextValue $ws.Range("G12") "10"

# Row 13
Set-TextValue $ws.Range("B13") "BitrueCoin"
Set-TextValue $ws.Range("C13") "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03612"
Set-TextValue $ws.Range("E13") "0.12%"
Set-TextValue $ws.Range("G13") "10"

# Row 14
Set-TextValue $ws.Range("B14") "BitMartToken"
Set-TextValue $ws.Range("C14") "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09859"
Set-TextValue $ws.Range("E14") "-0.48%"
Set-TextValue $ws.Range("G14") "10"

# Row 15
Set-TextValue $ws.Range("B15") "BitForexToken"
Set-TextValue $ws.Range("C15") "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001386"
Set-TextValue $ws.Range("E15") "-3.78%"
Set-TextValue $ws.Range("G15") "10"

# Row 16
Set-TextValue $ws.Range("B16") "TigerCash"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.005811"
Set-TextValue $ws.Range("E16") "1.36%"
Set-TextValue $ws.Range("G16") "10"

# Row 17
Set-TextValue $ws.Range("B17") "LEO"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.497"
Set-TextValue $ws.Range("E17") "0.84%"
Set-TextValue $ws.Range("G17") "10"

# Row 18
Set-TextValue $ws.Range("B18") "GateToken"
Set-TextValue $ws.Range("C18") "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D18") "4.039"
Set-TextValue $ws.Range("E18") "-2.69%"
Set-TextValue $ws.Range("G18") "10"

# Row 19
Set-TextValue $ws.Range("D19") "0.3429"
Set-TextValue $ws.Range("E19") "1.84%"
Set-TextValue $ws.Range("G19") "10"

# Row 20
Set-TextValue $ws.Range("E20") "-1.61%"
Set-TextValue $ws.Range("G20") "10"

# Row 21
Set-TextValue $ws.Range("D21") "5.049"
Set-TextValue $ws.Range("E21") "-1.76%"
Set-TextValue $ws.Range("G21") "10"

# Row 22
Set-TextValue $ws.Range("D22") "0.2471"
Set-TextValue $ws.Range("E22") "5.56%"
Set-TextValue $ws.Range("G22") "10"

# Row 23
Set-TextValue $ws.Range("D23") "0.04530"
Set-TextValue $ws.Range("G23") "10"

# Row 24
Set-TextValue $ws.Range("D24") "0.001219"
Set-TextValue $ws.Range("E24") "-2.21%"
Set-TextValue $ws.Range("G24") "10"

# Row 25
Set-TextValue $ws.Range("D25") "0.004798"
Set-TextValue $ws.Range("E25") "2.09%"
Set-TextValue $ws.Range("G25") "10"

# Row 26
Set-TextValue $ws.Range("D26") "0.0001254"
Set-TextValue $ws.Range("E26") "0.02%"
Set-TextValue $ws.Range("G26") "10"

# Row 27
Set-TextValue $ws.Range("D27") "0.0003007"
Set-TextValue $ws.Range("E27") "-33.34%"
Set-TextValue $ws.Range("G27") "10"

# Row 28
Set-TextValue $ws.Range("G28") "10"

# Row 29
Set-TextValue $ws.Range("G29") "10"

# Row 30
Set-TextValue $ws.Range("G30") "10"

# Row 31
Set-TextValue $ws.Range("G31") "10"

# Row 32
Set-TextValue $ws.Range("G32") "10"

# Row 33
Set-TextValue $ws.Range("G33") "10"

# Row 34
Set-TextValue $ws.Range("G34") "10"

# Row 35
Set-TextValue $ws.Range("G35") "10"

# Row 36
Set-TextValue $ws.Range("G36") "10"

# Row 37
Set-TextValue $ws.Range("G37") "10"

# Row 38
Set-TextValue $ws.Range("G38") "10"

# Row 39
Set-TextValue $ws.Range("D39") "0.01889"
Set-TextValue $ws.Range("G39") "10"

# Row 40
Set-TextValue $ws.Range("D40") "0.04681"
Set-TextValue $ws.Range("E40") "-4.00%"
Set-TextValue $ws.Range("G40") "10"

# Row 41
Set-TextValue $ws.Range("D41") "0.007553"
Set-TextValue $ws.Range("E41") "-3.21%"
Set-TextValue $ws.Range("G41") "10"

# Row 42
Set-TextValue $ws.Range("D42") "0.009709"
Set-TextValue $ws.Range("E42") "23.99%"
Set-TextValue $ws.Range("G42") "10"

# Row 43
Set-TextValue $ws.Range("D43") "0.1324"
Set-TextValue $ws.Range("E43") "-4.72%"
Set-TextValue $ws.Range("G43") "10"

# Row 44
Set-TextValue $ws.Range("D44") "0.002116"
Set-TextValue $ws.Range("E44") "-1.83%"
Set-TextValue $ws.Range("G44") "10"

# Row 45
Set-TextValue $ws.Range("D45") "0.01091"
Set-TextValue $ws.Range("E45") "-6.08%"
Set-TextValue $ws.Range("G45") "10"

# Row 46
Set-TextValue $ws.Range("D46") "0.00006201"
Set-TextValue $ws.Range("E46") "-4.73%"
Set-TextValue $ws.Range("G46") "10"

# Row 47
Set-TextValue $ws.Range("D47") "0.00000000752"
Set-TextValue $ws.Range("E47") "0.04%"
Set-TextValue $ws.Range("G47") "10"

# Row 48
Set-TextValue $ws.Range("E48") "79.90%"
Set-TextValue $ws.Range("G48") "10"

# Row 49
Set-TextValue $ws.Range("E49") "-21.92%"
Set-TextValue $ws.Range("G49") "10"

# Row 50
Set-TextValue $ws.Range("D50") "0.00002106"
Set-TextValue $ws.Range("E50") "0.04%"
Set-TextValue $ws.Range("G50") "10"

# Row 51
Set-TextValue $ws.Range("D51") "0.0002006"
Set-TextValue $ws.Range("E51") "0.04%"
Set-TextValue $ws.Range("G51") "10"
